$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - sheet1
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1776
$ws1.Range("F7").Value = 628
$ws1.Range("F8").Value = 324
$ws1.Range("F9").Value = 1712
$ws1.Range("F10").Value = 357
$ws1.Range("F11").Value = 1410
$ws1.Range("F13").Value = 331
$ws1.Range("F15").Value = 12730
$ws1.Range("F16").Value = 12750
$ws1.Range("F17").Value = 948
$ws1.Range("F18").Value = 742
$ws1.Range("F20").Value = 509
$ws1.Range("F21").Value = 50
$ws1.Range("F22").Value = 542
$ws1.Range("F23").Value = 1994
$ws1.Range("F26").Value = 239
$ws1.Range("F27").Value = 669

# Sheet "演出" (Performance) - sheet2
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 5
$ws2.Range("F5").Value = 77
$ws2.Range("F10").Value = 75

# Sheet "全部类型" (All Types) - sheet4
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1776
$ws4.Range("F7").Value = 5
$ws4.Range("F11").Value = 628
$ws4.Range("F13").Value = 324
$ws4.Range("F14").Value = 1712
$ws4.Range("F15").Value = 357
$ws4.Range("F16").Value = 1410
$ws4.Range("F18").Value = 331
$ws4.Range("F19").Value = 77
$ws4.Range("F21").Value = 12730
$ws4.Range("F22").Value = 12750
$ws4.Range("F23").Value = 948
$ws4.Range("F24").Value = 742
$ws4.Range("F26").Value = 509
$ws4.Range("F27").Value = 50
$ws4.Range("F28").Value = 542
$ws4.Range("F31").Value = 1995
$ws4.Range("F36").Value = 239
$ws4.Range("F37").Value = 669
$ws4.Range("F38").Value = 75
